# Updates cryptos list data (prices, volumes, and a couple of row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.916.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.46%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.764.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.66%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'624.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.76%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'165.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.13%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.763.17"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.65%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.19%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.56%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +2.84%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'6.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.41%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.10%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'35.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.60%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.400.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.53%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.743.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.36%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'68.900.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.50%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'17.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.32%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -1.20%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.25%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'466.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.73%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +1.03%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.704"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.93%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.0000145"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.61%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'82.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.41%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'12.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.13%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.78%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'Dai"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.08%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'RenderToken"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'9.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.40%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'3.914.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "'  +2.47%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +2.17%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'7.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.26%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'28.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'0.175"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +20.41%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.07%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.716.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.62%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'8.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.00%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +1.75%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.81%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.18%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.964"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.34%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.07%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D45").Value = "'153.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.01%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'43.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.59%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'OKB"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'46.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.91%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'TheGraph"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.294"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.05%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +3.46%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'8.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.11%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.50%  "
$ws.Range("E51").Style = "Normal"
